$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing Standard Deviation row (row 3) with refined values
$ws.Range("B3").Value = 0.004156683611941256
$ws.Range("C3").Value = 0.003026023945209601
$ws.Range("D3").Value = 0.002366190155367758

# Update existing Maximum row (row 4) with refined values
$ws.Range("B4").Value = 0.0168710555368855
$ws.Range("C4").Value = 0.01240980411715342
$ws.Range("D4").Value = 0.01259769366258634

# Add new Mean row (row 5)
$ws.Range("A5").Value = "Mean"
$ws.Range("B5").Value = 0.00570089522340696
$ws.Range("C5").Value = 0.004273305666607591
$ws.Range("D5").Value = 0.003546191850153032
